# NhapNoiThat.xlsx — "update pdf - 1"
#
# The source diff shows the header label in C1 ("nameFurniture") being
# blanked out to a single space, the body font switching from Calibri to
# Arial, and the active selection moving to E2. (The diff also contains a
# batch of Excel-session artifacts — xr:revisionPtr GUIDs, the x15ac
# absPath, bookViews window geometry, and sub-pixel column/row roundings
# from a font-substitution relayout — that are side effects of which
# machine/session last touched the file rather than deliberate edits, so
# they are not something a COM script can/should try to fabricate.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Clear the "nameFurniture" header text in C1, replacing it with a
#    single blank space (matches the shared-string edit in the diff).
$ws.Range("C1").Value = " "

# 2. Switch the data rows' font from Calibri to Arial. Row 1 (the header)
#    keeps its own bold font untouched, so restrict this to the data body
#    (A2:E11) to avoid clobbering the header's distinct style.
$ws.Range("A2:E11").Font.Name = "Arial"

# 3. Leave the selection on E2, matching the saved cursor position.
$ws.Range("E2").Select()
